# repull data, push all data, mean calculation
# Update column F (dSF) values for the data rows with newly pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -1
    4  = -5
    5  = 8
    6  = 2
    7  = 1
    8  = -3
    9  = 5
    10 = -1
    11 = -3
    12 = -3
    14 = -2
    15 = -6
    16 = -2
    17 = 5
    18 = 3
    19 = -4
    20 = -4
    21 = -1
    22 = 1
    23 = -1
    24 = -1
    25 = 3
    26 = -1
    27 = -1
    29 = -2
    30 = -2
    31 = 1
    33 = 4
    34 = -2
    35 = 1
    36 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
